$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()
$ws2.Range("A1").Value = "DeviceCount"
Write-Output "done"
